# Rename "Customer-Import" sheet to "Sold to Party-Import"
# (commit: Rename "Customer-Import" column to "Sold-to-party")

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Customer-Import").Name = "Sold to Party-Import"

# Browsing through the workbook leaves each visited sheet's selection
# reset back to A2, ending with the first sheet active again.
$wsDelivery = $wb.Worksheets.Item("Delivery-Point-Import")
$wsDelivery.Activate()
[void]$wsDelivery.Range("A2").Select()

$wsProduct = $wb.Worksheets.Item("Product-Import")
$wsProduct.Activate()
[void]$wsProduct.Range("A2").Select()

$wsOrg = $wb.Worksheets.Item("Organization-Import")
$wsOrg.Activate()
[void]$wsOrg.Range("A2").Select()
